# "updated with new sensor launch"
# Append three new sensor-delivery rows (37-39) to Sheet1 and refresh the
# sheet view (selection) to reflect where the user ended up working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the existing date-cell formatting (numFmtId 15 / style index 1)
# for the new "Date Delivered to Xin" cells.
$dateFormat = $ws.Range("F2").NumberFormat

# --- Row 37: new WS20-CXA sensor, delivered with PT15-6A5 below -----------
$ws.Range("A37").Value = "WS20-CXA"
$ws.Range("B37").Value = "65c3c5102e3cdb000d09150d"
$ws.Range("C37").Value = "65c3c51148bb6b000d333ba5"
$ws.Range("D37").Value = 40.307834999999997
$ws.Range("E37").Value = -104.777721
$ws.Range("F37").Value = 45482
$ws.Range("F37").NumberFormat = $dateFormat

# --- Row 38: re-delivery of the existing WS27-XLU sensor -------------------
$ws.Range("A38").Value = "WS27-XLU"
$ws.Range("B38").Value = "65d6457c7a715d000bf94dc0"
$ws.Range("C38").Value = "65d6457d7a715d000c7d068c"
$ws.Range("D38").Value = 39.962213920000003
$ws.Range("E38").Value = -102.29797979999999
$ws.Range("F38").Value = 45482
$ws.Range("F38").NumberFormat = $dateFormat

# --- Row 39: new PT15-6A5 sensor -------------------------------------------
$ws.Range("A39").Value = "PT15-6A5"
$ws.Range("B39").Value = "63fd2c8f4b3d79000b91720d"
$ws.Range("C39").Value = "6440333064e2e2000d87d60d"
$ws.Range("D39").Value = 40.498097999999999
$ws.Range("E39").Value = -104.59709599999999
$ws.Range("F39").Value = 45482
$ws.Range("F39").NumberFormat = $dateFormat

# --- "Delivered by" column filled in afterwards, per row -------------------
$ws.Range("G37").Value = "Terry Weber"
$ws.Range("G38").Value = "Stewart Norrish"
$ws.Range("G39").Value = "Terry Weber"

# --- Scroll/selection state reflecting where the user left off -------------
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
[void]$ws.Range("B43").Select()
